$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> [C, D, E] new values (Total = B column is untouched)
$updates = @{
    2  = @(81, 81, 81)
    3  = @(154, 135, 122)
    4  = @(235, 183, 161)
    5  = @(31, 31, 29)
    7  = @(100, 81, 70)
    8  = @(77, 77, 64)
    9  = @(82, 76, 69)
    10 = @(21, 21, 19)
    11 = @(138, 124, 114)
    12 = @(239, 179, 145)
    13 = @(170, 120, 98)
    14 = @(45, 43, 39)
    15 = @(22, 22, 21)
    16 = @(14, 14, 14)
    30 = @(21, 21, 21)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]  # Column C
    $ws.Cells.Item($row, 4).Value = $vals[1]  # Column D
    $ws.Cells.Item($row, 5).Value = $vals[2]  # Column E
}

$wb.Save()
